$wb = $excel.ActiveWorkbook

# --- Add Sheet2 after Sheet1 ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Sheet2"

# Matching column widths used on Sheet1
$ws.Columns.Item(2).ColumnWidth = 18.625
$ws.Columns.Item(3).ColumnWidth = 18.25

# Matching page setup used on Sheet1
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Row 1: title (merged A1:C1, centered) ---
$ws.Range("A1").Value = "title"
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A1:C1").VerticalAlignment = -4108
$ws.Range("A1:C1").Merge()

# --- Row 2: subtitle1 (merged A2:B2, centered) / subtitle2 (C2, centered) ---
$ws.Range("A2").Value = "subtitle1"
$ws.Range("C2").Value = "subtitle2"
$ws.Range("A2:B2").HorizontalAlignment = -4108
$ws.Range("A2:B2").VerticalAlignment = -4108
$ws.Range("A2:B2").Merge()
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").VerticalAlignment = -4108

# --- Row 3: header row ---
$ws.Range("A3").Value = "id"
$ws.Range("B3").Value = "name"
$ws.Range("C3").Value = "email"

# --- Rows 4-13: sample/option data ---
# ids (column A) for rows 4..13 -> 1..10
$ids = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$r = 4
foreach ($id in $ids) {
    $ws.Cells.Item($r, 1).Value = $id
    $r = $r + 1
}

# names (column B) - filled top to bottom first, row 7 (id 4) has none
$names = @("1_option_name_xlsx", "2_option_name_xlsx", "3_option_name_xlsx", $null, "5_option_name_xlsx", "6_option_name_xlsx", "7_option_name_xlsx", "8_option_name_xlsx", "9_option_name_xlsx", "10_option_name_xlsx")
$r = 4
foreach ($name in $names) {
    if ($name -ne $null) { $ws.Cells.Item($r, 2).Value = $name }
    $r = $r + 1
}

# emails (column C) - filled top to bottom next, row 10 (id 7) has none
$emails = @("1_option_xlsx@email.com", "2_option_xlsx@email.com", "3_option_xlsx@email.com", "4_option_xlsx@email.com", "5_option_xlsx@email.com", "6_option_xlsx@email.com", $null, "8_option_xlsx@email.com", "9_option_xlsx@email.com", "10_option_xlsx@email.com")
$r = 4
foreach ($email in $emails) {
    if ($email -ne $null) { $ws.Cells.Item($r, 3).Value = $email }
    $r = $r + 1
}

# --- Selection on the new sheet ---
$ws.Range("C10").Select()

# --- Sheet1 is no longer the tab shown on open; Sheet2 (index 2) is active ---
$ws.Activate()
